$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 402:403, pushing existing rows 402-432 down to 404-434
$ws.Rows("402:403").Insert()

# New row 402 - "Primera" quality, week of 45013
$ws.Range("A402").Value = 4
$ws.Range("B402").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C402").Value = "Los Lagos"
$ws.Range("D402").Value = 45013
$ws.Range("E402").Value = 10
$ws.Range("F402").Value = 100112017
$ws.Range("G402").Value = "Apio"
$ws.Range("H402").Value = "Americana (o)"
$ws.Range("I402").Value = "Primera"
$ws.Range("J402").Value = 25
$ws.Range("K402").Value = 11000
$ws.Range("L402").Value = 11000
$ws.Range("M402").Value = 11000
$ws.Range("N402").Value = "`$/docena de matas"
$ws.Range("O402").Value = "Región de Coquimbo"
$ws.Range("P402").Value = 1833
$ws.Range("Q402").Value = 6
$ws.Range("R402").Value = "Hortaliza"

# New row 403 - "Segunda" quality, week of 45013
$ws.Range("A403").Value = 4
$ws.Range("B403").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C403").Value = "Los Lagos"
$ws.Range("D403").Value = 45013
$ws.Range("E403").Value = 10
$ws.Range("F403").Value = 100112017
$ws.Range("G403").Value = "Apio"
$ws.Range("H403").Value = "Americana (o)"
$ws.Range("I403").Value = "Segunda"
$ws.Range("J403").Value = 25
$ws.Range("K403").Value = 10000
$ws.Range("L403").Value = 10000
$ws.Range("M403").Value = 10000
$ws.Range("N403").Value = "`$/docena de matas"
$ws.Range("O403").Value = "Región de Coquimbo"
$ws.Range("P403").Value = 1667
$ws.Range("Q403").Value = 6
$ws.Range("R403").Value = "Hortaliza"
